$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (paragraph 1, Heading1 "Play Angel's Touch Free Slot Game | Review 2021").
#    The new paragraph has no paragraph style (like the other body paragraphs),
#    starts with an empty run, then a bold "Meta description" run, then a
#    normal run with the rest of the sentence.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Looking to play Angel' + [char]39 + 's Touch for free? Check out our review of this celestial-themed slot game and find out what we like and don' + [char]39 + 't like about it.</w:t></w:r>' +
  '</w:p>'

$null = $metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the old bold "Play Angel's Touch Free Slot Game | Review 2021"
#    paragraph that used to sit right before the italic "Looking to play..."
#    paragraph near the end of the document. (The Heading1 title paragraph at
#    the very top has the identical text, so match on exact paragraph text
#    and skip the Heading1 occurrence.)
# ---------------------------------------------------------------------------
$titleText = "Play Angel's Touch Free Slot Game | Review 2021"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $pr = $p.Range
    $txt = $pr.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $titleText -and $p.Style.NameLocal -ne "Heading 1") {
        $pr.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the remaining italic paragraph with the new
#    image-prompt text, keeping its italic run formatting intact.
# ---------------------------------------------------------------------------
$oldExact = "Looking to play Angel's Touch for free? Check out our review of this celestial-themed slot game and find out what we like and don't like about it."
$newText = "Create an eye-catching feature image for Angel's Touch that will catch the attention of online slot game enthusiasts. The image should be in a fun cartoon style and showcase a happy Maya warrior with glasses. The background should be bright and captivating with an angelic theme. The image must include the game's logo and the title in a bold and easy-to-read font. Make sure the Maya warrior is engaged in an exciting activity, like spinning the reels or celebrating a big win while surrounded by angel wings. The overall look and feel of the image should be fun, exciting, and attention-grabbing to encourage players to try their luck with Angel's Touch."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pr = $p.Range
    $txt = $pr.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $oldExact) {
        $target = $d.Range($pr.Start, $pr.End - 1)
        $target.Text = $newText
        break
    }
}
